$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 1991
$ws.Range("E2").Value = 324
$ws.Range("F2").Value = 324
$ws.Range("G2").Value = 336
$ws.Range("H2").Value = 238
$ws.Range("I2").Value = 238
$ws.Range("K2").Value = 2194
$ws.Range("L2").Value = 428
$ws.Range("M2").Value = 1766
$ws.Range("N2").Value = 1766
$ws.Range("P2").Value = 96
$ws.Range("Q2").Value = 385
$ws.Range("R2").Value = -368
$ws.Range("S2").Value = -73
$ws.Range("T2").Value = 138
$ws.Range("U2").Value = 247
$ws.Range("V2").Value = 1
$ws.Range("W2").Value = 16.29
$ws.Range("X2").Value = 11.97
$ws.Range("Y2").Value = 14.32
$ws.Range("Z2").Value = 11.29
$ws.Range("AA2").Value = 24.23
$ws.Range("AB2").Value = 1719.67
$ws.Range("AC2").Value = 1241
$ws.Range("AD2").Value = 8.859999999999999
$ws.Range("AE2").Value = 9198
$ws.Range("AF2").Value = 1.2
$ws.Range("AG2").Value = 200
$ws.Range("AH2").Value = 1.82
$ws.Range("AI2").Value = 16.11
$ws.Range("AJ2").Value = 19200000
$ws.Range("J2").ClearContents()
$ws.Range("O2").ClearContents()

# Row 3
$ws.Range("D3").Value = 2027
$ws.Range("E3").Value = 297
$ws.Range("F3").Value = 297
$ws.Range("G3").Value = 328
$ws.Range("H3").Value = 248
$ws.Range("I3").Value = 248
$ws.Range("K3").Value = 2460
$ws.Range("L3").Value = 482
$ws.Range("M3").Value = 1978
$ws.Range("N3").Value = 1978
$ws.Range("P3").Value = 96
$ws.Range("Q3").Value = 346
$ws.Range("R3").Value = -253
$ws.Range("S3").Value = -38
$ws.Range("T3").Value = 112
$ws.Range("U3").Value = 234
$ws.Range("V3").Value = 1
$ws.Range("W3").Value = 14.67
$ws.Range("X3").Value = 12.22
$ws.Range("Y3").Value = 13.24
$ws.Range("Z3").Value = 10.65
$ws.Range("AA3").Value = 24.38
$ws.Range("AB3").Value = 1928
$ws.Range("AC3").Value = 1291
$ws.Range("AD3").Value = 6.58
$ws.Range("AE3").Value = 10301
$ws.Range("AF3").Value = 0.82
$ws.Range("AG3").Value = 200
$ws.Range("AH3").Value = 2.36
$ws.Range("AI3").Value = 15.5
$ws.Range("AJ3").Value = 19200000
$ws.Range("J3").ClearContents()
$ws.Range("O3").ClearContents()

# Row 4
$ws.Range("D4").Value = 1975
$ws.Range("E4").Value = 280
$ws.Range("F4").Value = 280
$ws.Range("G4").Value = 314
$ws.Range("H4").Value = 226
$ws.Range("I4").Value = 226
$ws.Range("K4").Value = 2637
$ws.Range("L4").Value = 496
$ws.Range("M4").Value = 2141
$ws.Range("N4").Value = 2141
$ws.Range("P4").Value = 96
$ws.Range("Q4").Value = 387
$ws.Range("R4").Value = -223
$ws.Range("S4").Value = -38
$ws.Range("T4").Value = 90
$ws.Range("U4").Value = 297
$ws.Range("V4").Value = 1
$ws.Range("W4").Value = 14.2
$ws.Range("X4").Value = 11.42
$ws.Range("Y4").Value = 10.95
$ws.Range("Z4").Value = 8.85
$ws.Range("AA4").Value = 23.16
$ws.Range("AB4").Value = 2126.97
$ws.Range("AC4").Value = 1175
$ws.Range("AD4").Value = 7.08
$ws.Range("AE4").Value = 11153
$ws.Range("AF4").Value = 0.75
$ws.Range("AG4").Value = 190
$ws.Range("AH4").Value = 2.28
$ws.Range("AI4").Value = 16.17
$ws.Range("AJ4").Value = 19200000
$ws.Range("J4").ClearContents()
$ws.Range("O4").ClearContents()

# Row 5
$ws.Range("D5").Value = 1835
$ws.Range("E5").Value = 167
$ws.Range("F5").Value = 167
$ws.Range("G5").Value = 191
$ws.Range("H5").Value = 133
$ws.Range("I5").Value = 133
$ws.Range("K5").Value = 2666
$ws.Range("L5").Value = 470
$ws.Range("M5").Value = 2196
$ws.Range("N5").Value = 2196
$ws.Range("P5").Value = 96
$ws.Range("Q5").Value = 87
$ws.Range("R5").Value = 88
$ws.Range("S5").Value = -36
$ws.Range("T5").Value = 155
$ws.Range("U5").Value = -68
$ws.Range("V5").Value = 1
$ws.Range("W5").Value = 9.1
$ws.Range("X5").Value = 7.23
$ws.Range("Y5").Value = 6.12
$ws.Range("Z5").Value = 5
$ws.Range("AA5").Value = 21.39
$ws.Range("AB5").Value = 2228.9
$ws.Range("AC5").Value = 691
$ws.Range("AD5").Value = 10.85
$ws.Range("AE5").Value = 11439
$ws.Range("AF5").Value = 0.66
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 2
$ws.Range("AI5").Value = 21.71
$ws.Range("AJ5").Value = 19200000
$ws.Range("J5").ClearContents()
$ws.Range("O5").ClearContents()

# Row 6
$ws.Range("D6").Value = 1676
$ws.Range("E6").Value = 95
$ws.Range("F6").Value = 95
$ws.Range("G6").Value = 100
$ws.Range("H6").Value = 74
$ws.Range("I6").Value = 74
$ws.Range("K6").Value = 2671
$ws.Range("L6").Value = 447
$ws.Range("M6").Value = 2224
$ws.Range("N6").Value = 2224
$ws.Range("P6").Value = 96
$ws.Range("Q6").Value = 319
$ws.Range("R6").Value = -31
$ws.Range("S6").Value = -29
$ws.Range("T6").Value = 357
$ws.Range("U6").Value = -38
$ws.Range("V6").Value = 1
$ws.Range("W6").Value = 5.7
$ws.Range("X6").Value = 4.41
$ws.Range("Y6").Value = 3.35
$ws.Range("Z6").Value = 2.77
$ws.Range("AA6").Value = 20.1
$ws.Range("AB6").Value = 2270.7
$ws.Range("AC6").Value = 385
$ws.Range("AD6").Value = 15.3
$ws.Range("AE6").Value = 11581
$ws.Range("AF6").Value = 0.51
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 2.55
$ws.Range("AI6").Value = 38.96
$ws.Range("AJ6").Value = 19200000

# Row 7
$ws.Range("D7").Value = 1548
$ws.Range("E7").Value = -37
$ws.Range("G7").Value = 23
$ws.Range("H7").Value = 24
$ws.Range("I7").Value = 24
$ws.Range("K7").Value = 2633
$ws.Range("L7").Value = 415
$ws.Range("M7").Value = 2219
$ws.Range("N7").Value = 2219
$ws.Range("P7").Value = 96
$ws.Range("Q7").Value = 118
$ws.Range("R7").Value = -110
$ws.Range("S7").Value = -29
$ws.Range("T7").Value = 150
$ws.Range("U7").Value = -32
$ws.Range("W7").Value = -2.39
$ws.Range("X7").Value = 1.55
$ws.Range("Y7").Value = 1.08
$ws.Range("Z7").Value = 0.91
$ws.Range("AA7").Value = 18.7
$ws.Range("AC7").Value = 125
$ws.Range("AD7").Value = 36.4
$ws.Range("AE7").Value = 11557
$ws.Range("AF7").Value = 0.39
$ws.Range("AG7").Value = 150
$ws.Range("AH7").Value = 3.3
$ws.Range("AI7").Value = 120

# Row 8
$ws.Range("D8").Value = 1598
$ws.Range("E8").Value = 13
$ws.Range("G8").Value = 55
$ws.Range("H8").Value = 41
$ws.Range("I8").Value = 41
$ws.Range("K8").Value = 2659
$ws.Range("L8").Value = 427
$ws.Range("M8").Value = 2231
$ws.Range("N8").Value = 2231
$ws.Range("P8").Value = 96
$ws.Range("Q8").Value = 148
$ws.Range("R8").Value = -84
$ws.Range("S8").Value = -29
$ws.Range("T8").Value = 80
$ws.Range("U8").Value = 68
$ws.Range("W8").Value = 0.8100000000000001
$ws.Range("X8").Value = 2.57
$ws.Range("Y8").Value = 1.84
$ws.Range("Z8").Value = 1.55
$ws.Range("AA8").Value = 19.14
$ws.Range("AC8").Value = 214
$ws.Range("AD8").Value = 21.31
$ws.Range("AE8").Value = 11620
$ws.Range("AF8").Value = 0.39
$ws.Range("AG8").Value = 150
$ws.Range("AH8").Value = 3.3
$ws.Range("AI8").Value = 70.23999999999999

# Row 9
$ws.Range("D9").Value = 1678
$ws.Range("E9").Value = 27
$ws.Range("G9").Value = 79
$ws.Range("H9").Value = 59
$ws.Range("I9").Value = 59
$ws.Range("K9").Value = 2709
$ws.Range("L9").Value = 448
$ws.Range("M9").Value = 2261
$ws.Range("N9").Value = 2261
$ws.Range("P9").Value = 96
$ws.Range("Q9").Value = 148
$ws.Range("R9").Value = -81
$ws.Range("S9").Value = -29
$ws.Range("T9").Value = 80
$ws.Range("U9").Value = 68
$ws.Range("W9").Value = 1.61
$ws.Range("X9").Value = 3.52
$ws.Range("Y9").Value = 2.63
$ws.Range("Z9").Value = 2.2
$ws.Range("AA9").Value = 19.81
$ws.Range("AC9").Value = 307
$ws.Range("AD9").Value = 14.81
$ws.Range("AE9").Value = 11776
$ws.Range("AF9").Value = 0.39
$ws.Range("AG9").Value = 150
$ws.Range("AH9").Value = 3.3
$ws.Range("AI9").Value = 48.81
